$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.327.21'
$ws.Range("E2").Value = '  -2.52%  '

$ws.Range("D3").Value = '2.199.97'
$ws.Range("E3").Value = '  -6.93%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.86%  '

$ws.Range("E7").Value = '  -4.42%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("E9").Value = '  -4.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0770'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.94'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.33'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -10.24%  '

$ws.Range("E13").Value = '  -2.44%  '

$ws.Range("D14").Value = '2.544.11'
$ws.Range("E14").Value = '  -7.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.25'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.01%  '

$ws.Range("E16").Value = '  -5.77%  '

$ws.Range("D17").Value = '2.197.71'
$ws.Range("E17").Value = '  -6.81%  '

$ws.Range("E18").Value = '  -5.77%  '

$ws.Range("D19").Value = '39.191.37'
$ws.Range("E19").Value = '  -2.94%  '

$ws.Range("E20").Value = '  -3.79%  '

$ws.Range("E21").Value = '  -6.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '225.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.92%  '

$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("E26").Value = '  -6.64%  '

$ws.Range("E27").Value = '  -0.45%  '

$ws.Range("E28").Value = '  -4.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.15%  '

$ws.Range("E33").Value = '  -0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.82%  '

$ws.Range("E35").Value = '  -3.65%  '

$ws.Range("E37").Value = '  -3.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0964'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.08'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.50%  '

$ws.Range("E40").Value = '  -5.91%  '

$ws.Range("E41").Value = '  -3.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.61'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.85%  '

$ws.Range("D43").Value = '1.892.16'
$ws.Range("E43").Value = '  -3.63%  '

$ws.Range("E44").Value = '  -12.18%  '

$ws.Range("E45").Value = '  -3.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.42%  '

$ws.Range("E48").Value = '  -2.90%  '

$ws.Range("D49").Value = '2.415.72'
$ws.Range("E49").Value = '  -7.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.53%  '
